$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 41.364051818847656
$ws.Range("B3").Value = 37.248386383056641
$ws.Range("B4").Value = 18.59356689453125
$ws.Range("B5").Value = 7.5408868789672852
$ws.Range("B6").Value = 6.0619406700134277
$ws.Range("B7").Value = 14.847681045532227
$ws.Range("B8").Value = 20.680257797241211
$ws.Range("B9").Value = 7.4842095375061035
$ws.Range("B10").Value = 32.345989227294922
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
